# Updated cryptos list on Mon Sep 18 08:11:01 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.862.58"
$ws.Range("E2").Value = "  +0.28%  "
# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.38"
$ws.Range("E3").Value = "  -0.19%  "
# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.49%  "
# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.17"
$ws.Range("E5").Value = "  +0.45%  "
# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.496"
$ws.Range("E6").Value = "  -0.99%  "
# Row 7 - USDC
$ws.Range("E7").Value = "  -0.47%  "
# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.57%  "
# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.03%  "
# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.24"
$ws.Range("E10").Value = "  +0.62%  "
# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.15%  "
# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.74"
$ws.Range("E12").Value = "  -0.08%  "
# Row 13 - WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.649.30"
$ws.Range("E13").Value = "  +0.01%  "
# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.23%  "
# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.07%  "
# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.12"
$ws.Range("E16").Value = "  +0.99%  "
# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.866.80"
$ws.Range("E17").Value = "  +0.23%  "
# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.73%  "
# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.21"
$ws.Range("E19").Value = "  +0.71%  "
# Row 20 - Dai
$ws.Range("E20").Value = "  -0.47%  "
# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.04%  "
# Row 22 - Chainlink
$ws.Range("E22").Value = "  +4.89%  "
# Row 23 - Toncoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.38"
$ws.Range("E23").Value = "  -4.00%  "
# Row 24 - Avalanche
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  -1.44%  "
# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.37"
$ws.Range("E25").Value = "  +1.29%  "
# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.41%  "
# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.54%  "
# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("E28").Value = "  +0.84%  "
# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +0.29%  "
# Row 30 - Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -0.06%  "
# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.77%  "
# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +0.85%  "
# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.20%  "
# Row 34 - Maker
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.282.26"
# Row 35 - LidoDAOToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +1.08%  "
# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.08%  "
# Row 37 - VeChain
$ws.Range("E37").Value = "  -1.32%  "
# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -0.51%  "
# Row 39 - ARBITRUM
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  -0.05%  "
# Row 40 - PaxDollar
$ws.Range("E40").Value = "  -0.40%  "
# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -0.60%  "
# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.33"
$ws.Range("E42").Value = "  +0.12%  "
# Row 43 - RocketPoolETH
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.782.95"
# Row 44 - MXToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.10"
$ws.Range("E44").Value = "  -6.06%  "
# Row 45 - Quant
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.45"
$ws.Range("E45").Value = "  +1.07%  "
# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.04"
$ws.Range("E46").Value = "  -1.15%  "
# Row 47 - RenderToken
$ws.Range("E47").Value = "  -1.09%  "
# Row 48 - Cronos
$ws.Range("E48").Value = "  -1.76%  "

# Rows 49/50 - Algorand and EnergySwap swap places (with updated values)
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.59"
$ws.Range("E49").Value = "  -0.90%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0968"
$ws.Range("E50").Value = "  -0.54%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.26%  "
